$wb = $excel.ActiveWorkbook

# ---- constants describing the change (old GUID -> new GUID, hash change, new timestamps) ----
$oldGuid = "e4839f09-4783-4493-a9ef-43051d9ac33e"
$newGuid = "f7c5f93a-6302-4bf8-80fe-193effac2677"
$oldHash = "d24940cb3142328bcd4ea376c611b2164c726f53"
$newHash = "f6c6d10f2c2428ea0436459a440e31bdd6480fc0"

$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c79eb27e5c6a70b67ec7d8f1dcba03b5c941118f/e2e/$oldGuid.md"

# Target stored column width for col A is 39.3653477260045 characters.
# ColumnWidth is quantized by Excel to whole pixels (steps of 1/6 of a
# character here), so 38.5 is the closest input that reproduces the
# nearest achievable stored width (~39.333333333333336).
$newColAWidth = 38.5

# =========================================================================
# Overview sheet
# =========================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-19 19:05:49"

# Re-create hyperlink on B2 with the new display text (same target address)
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkAddress, [Type]::Missing, [Type]::Missing, "e2e\$newGuid.md") | Out-Null

$wsOverview.Columns.Item(1).ColumnWidth = $newColAWidth

# =========================================================================
# zh-cn sheet
# =========================================================================
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-19 19:05:45"

# Re-create hyperlink on A2 with the new display text (same target address)
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkAddress, [Type]::Missing, [Type]::Missing, "$newGuid.md") | Out-Null

$wsZhCn.Columns.Item(1).ColumnWidth = $newColAWidth

# =========================================================================
# de-de sheet
# =========================================================================
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-19 19:05:49"

# Re-create hyperlink on A2 with the new display text (same target address)
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkAddress, [Type]::Missing, [Type]::Missing, "$newGuid.md") | Out-Null

$wsDeDe.Columns.Item(1).ColumnWidth = $newColAWidth
